$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Fontana di Trevi photo URL (row 4, column F)
$ws.Range("F4").Value = "https://news.artnet.com/app/news-upload/2015/07/Fontana-di-Trevi.png"

# Add new row 6 (Reggio di Calabria) to the table and map.
# Copy the plain row formatting (row 3) down into row 6 first so the
# new row picks up the same style as the other non-hyperlinked rows.
$ws.Range("A3:I3").Copy()
$ws.Range("A6:I6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A6").Value = "Reggio di Calabria"
$ws.Range("B6").Value = "http://www.classicult.it/wp-content/uploads/2019/04/PHOTO-2018-03-16-21-23-20-1024x683.jpg"
$ws.Range("C6").Value = "The capital of the Calabrese Province. Here ‘L’Arenna dello Stretto’ in which you can find the statue of ‘Athena Promachos’ that protects the city."
$ws.Range("D6").Value = "rabarama.JPG"
$ws.Range("E6").Value = "This modern art is in the main road of this beautiful city. It was designed by the famous Italian artist Paola Rabarama."
$ws.Range("F6").Value = "https://upload.wikimedia.org/wikipedia/commons/d/d1/Il_prospetto_principale_del_duomo.jpg"
$ws.Range("G6").Value = "Reggio di Calabria’s Cathedral"
$ws.Range("H6").Value = 38.106266
$ws.Range("I6").Value = 15.647941
